$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 173.95653
$ws.Range("I33").Value = 181
$ws.Range("K33").Value = 181
$ws.Range("M33").Value = 48
# Row 40
$ws.Range("H40").Value = 1681.4166
$ws.Range("I40").Value = 1522.125
$ws.Range("K40").Value = 1522.125
$ws.Range("M40").Value = -1347.125
# Row 109
$ws.Range("H109").Value = 55286.6
$ws.Range("J109").Value = 55286.6
$ws.Range("L109").Value = 55286.6
$ws.Range("N109").Value = -58060.6
# Row 132
$ws.Range("H132").Value = 1594.1111
$ws.Range("I132").Value = 1431.5714
$ws.Range("J132").Value = 2163
$ws.Range("K132").Value = 4294.7142
$ws.Range("L132").Value = 6489
$ws.Range("M132").Value = -1764.7142
$ws.Range("N132").Value = -11549

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 13320.066
$ws.Range("I32").Value = 13537.812
$ws.Range("K32").Value = 13537.812
$ws.Range("M32").Value = -13250.812
# Row 45
$ws.Range("H45").Value = 1777.3334
$ws.Range("I45").Value = 1772.8
$ws.Range("K45").Value = 1772.8
$ws.Range("M45").Value = -1395.8
# Row 122
$ws.Range("H122").Value = 8186.3887
$ws.Range("I122").Value = 8620.5
$ws.Range("J122").Value = 6667
$ws.Range("K122").Value = 25861.5
$ws.Range("L122").Value = 20001
$ws.Range("M122").Value = -23411.5
$ws.Range("N122").Value = -24901
# Row 132
$ws.Range("H132").Value = 4674.0225
$ws.Range("I132").Value = 4742.147
$ws.Range("J132").Value = 4463.4546
$ws.Range("K132").Value = 14226.441
$ws.Range("L132").Value = 13390.3638
$ws.Range("M132").Value = -11696.441
$ws.Range("N132").Value = -18450.3638

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 107
$ws.Range("H107").Value = 23682.68
$ws.Range("I107").Value = 35013.625
$ws.Range("J107").Value = 3538.7778
$ws.Range("K107").Value = 35013.625
$ws.Range("L107").Value = 3538.7778
$ws.Range("M107").Value = -33093.625
$ws.Range("N107").Value = -7378.7778
# Row 134
$ws.Range("H134").Value = 3532.5789
$ws.Range("I134").Value = 3124.5386
$ws.Range("J134").Value = 4416.6665
$ws.Range("K134").Value = 9373.6158
$ws.Range("L134").Value = 13249.9995
$ws.Range("M134").Value = -6838.6158
$ws.Range("N134").Value = -18319.9995

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1719.1687
$ws.Range("I31").Value = 2085.3333
$ws.Range("K31").Value = 2085.3333
$ws.Range("M31").Value = -1790.3333
# Row 34
$ws.Range("H34").Value = 1719.1687
$ws.Range("I34").Value = 2085.3333
$ws.Range("K34").Value = 2085.3333
$ws.Range("M34").Value = -1883.3333
# Row 94
$ws.Range("H94").Value = 1132.381
$ws.Range("I94").Value = 661.3333
$ws.Range("J94").Value = 1485.6666
$ws.Range("K94").Value = 661.3333
$ws.Range("L94").Value = 1485.6666
$ws.Range("M94").Value = -210.3333
$ws.Range("N94").Value = -2387.6666
# Row 134
$ws.Range("H134").Value = 1359.9354
$ws.Range("I134").Value = 1157.2069
$ws.Range("J134").Value = 4299.5
$ws.Range("K134").Value = 3471.620699999999
$ws.Range("L134").Value = 12898.5
$ws.Range("M134").Value = -936.6206999999995
$ws.Range("N134").Value = -17968.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 791.2857
$ws.Range("J5").Value = 405.1579
$ws.Range("L5").Value = 1215.4737
$ws.Range("N5").Value = -1439.4737
# Row 64
$ws.Range("H64").Value = 3782
$ws.Range("I64").Value = 2841.6667
$ws.Range("J64").Value = 5035.778
$ws.Range("K64").Value = 8525.000100000001
$ws.Range("L64").Value = 15107.334
$ws.Range("M64").Value = -8255.000100000001
$ws.Range("N64").Value = -15647.334
# Row 67
$ws.Range("H67").Value = 3782
$ws.Range("I67").Value = 2841.6667
$ws.Range("J67").Value = 5035.778
$ws.Range("K67").Value = 8525.000100000001
$ws.Range("L67").Value = 15107.334
$ws.Range("M67").Value = -7589.000100000001
$ws.Range("N67").Value = -16979.334
# Row 68
$ws.Range("H68").Value = 916.29
$ws.Range("J68").Value = 1281.4736
$ws.Range("L68").Value = 3844.4208
$ws.Range("N68").Value = -5466.4208
# Row 71
$ws.Range("H71").Value = 916.29
$ws.Range("J71").Value = 1281.4736
$ws.Range("L71").Value = 11533.2624
$ws.Range("N71").Value = -19645.2624
# Row 107
$ws.Range("H107").Value = 1223.5692
$ws.Range("I107").Value = 1180.1538
$ws.Range("J107").Value = 1288.6923
$ws.Range("K107").Value = 3540.4614
$ws.Range("L107").Value = 3866.0769
$ws.Range("M107").Value = -1620.4614
$ws.Range("N107").Value = -7706.0769
# Row 117
$ws.Range("H117").Value = 373.55554
$ws.Range("J117").Value = 260.5
$ws.Range("L117").Value = 781.5
$ws.Range("N117").Value = -7665.5
# Row 120
$ws.Range("H120").Value = 11088.333
$ws.Range("J120").Value = 19000
$ws.Range("L120").Value = 57000
$ws.Range("N120").Value = -66676
# Row 135
$ws.Range("H135").Value = 791.2857
$ws.Range("J135").Value = 405.1579
$ws.Range("L135").Value = 3646.4211
$ws.Range("N135").Value = -8716.4211
# Row 140
$ws.Range("H140").Value = 1765.9062
$ws.Range("I140").Value = 880.95
$ws.Range("J140").Value = 3240.8333
$ws.Range("K140").Value = 2642.85
$ws.Range("L140").Value = 9722.499899999999
$ws.Range("M140").Value = 2537.15
$ws.Range("N140").Value = -20082.4999

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 25
$ws.Range("H25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").ClearContents()
# Row 57
$ws.Range("H57").Value = 28533.334
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 28533.334
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 28533.334
$ws.Range("M57").ClearContents()
$ws.Range("N57").Value = -30173.334
# Row 70
$ws.Range("H70").Value = 301385.25
$ws.Range("I70").Value = 407608.16
$ws.Range("K70").Value = 407608.16
$ws.Range("M70").Value = -407338.16
# Row 73
$ws.Range("H73").Value = 301385.25
$ws.Range("I73").Value = 407608.16
$ws.Range("K73").Value = 407608.16
$ws.Range("M73").Value = -406672.16
# Row 122
$ws.Range("H122").Value = 2701
$ws.Range("I122").Value = 2646.7273
$ws.Range("J122").Value = 2999.5
$ws.Range("K122").Value = 7940.1819
$ws.Range("L122").Value = 8998.5
$ws.Range("M122").Value = -5490.1819
$ws.Range("N122").Value = -13898.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 4
$ws.Range("H4").Value = 26000
$ws.Range("I4").Value = 26000
$ws.Range("K4").Value = 26000
$ws.Range("M4").Value = -25887
# Row 28
$ws.Range("H28").Value = 26000
$ws.Range("I28").Value = 26000
$ws.Range("K28").Value = 26000
$ws.Range("M28").Value = -25768
# Row 37
$ws.Range("H37").Value = 26000
$ws.Range("I37").Value = 26000
$ws.Range("K37").Value = 26000
$ws.Range("M37").Value = -25893
# Row 46
$ws.Range("H46").Value = 1172.1538
$ws.Range("I46").Value = 1022.7143
$ws.Range("J46").Value = 1346.5
$ws.Range("K46").Value = 1022.7143
$ws.Range("L46").Value = 1346.5
$ws.Range("M46").Value = -834.7143
$ws.Range("N46").Value = -1722.5
# Row 122
$ws.Range("H122").Value = 40915268
$ws.Range("I122").Value = 83337320
$ws.Range("J122").Value = 25006998
$ws.Range("K122").Value = 250011960
$ws.Range("L122").Value = 75020994
$ws.Range("M122").Value = -250009510
$ws.Range("N122").Value = -75025894
# Row 127
$ws.Range("H127").Value = 70355
$ws.Range("J127").Value = 70355
$ws.Range("L127").Value = 70355
$ws.Range("N127").Value = -80275

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 39
$ws.Range("H39").Value = 16060
$ws.Range("I39").Value = 100
$ws.Range("K39").Value = 100
$ws.Range("M39").Value = 313
# Row 42
$ws.Range("H42").Value = 46681.332
# Row 43
$ws.Range("H43").Value = 50000
$ws.Range("I43").Value = 50000
$ws.Range("K43").Value = 50000
$ws.Range("M43").Value = -49851
# Row 100
$ws.Range("H100").Value = 31300.2
$ws.Range("I100").Value = 100000.664
$ws.Range("J100").Value = 1857.1428
$ws.Range("K100").Value = 200001.328
$ws.Range("L100").Value = 3714.2856
$ws.Range("M100").Value = -199460.328
$ws.Range("N100").Value = -4796.2856
# Row 112
$ws.Range("H112").Value = 265000
$ws.Range("J112").Value = 265000
$ws.Range("L112").Value = 265000
$ws.Range("N112").Value = -267954
# Row 122
$ws.Range("H122").Value = 104168560
$ws.Range("I122").Value = 138890670
$ws.Range("J122").Value = 2218.6667
$ws.Range("K122").Value = 416672010
$ws.Range("L122").Value = 6656.000100000001
$ws.Range("M122").Value = -416669560
$ws.Range("N122").Value = -11556.0001
